# Apply odds updates to the FlashScore "Jogos da Semana" worksheet.
# The workbook contains a single worksheet (Sheet1) holding match odds;
# this script updates the specific cells that changed between the
# previous export and the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (Atletico-MG vs Gremio)
$ws.Range("G2").Value  = 1.95
$ws.Range("I2").Value  = 4
$ws.Range("J2").Value  = 2.63
$ws.Range("Z2").Value  = 17
$ws.Range("AC2").Value = 8.5
$ws.Range("AI2").Value = 19
$ws.Range("AQ2").Value = 41
$ws.Range("AW2").Value = 5.5

# Row 6 (Atl. Nacional vs Dep. Cali)
$ws.Range("Q6").Value  = 2.15
$ws.Range("R6").Value  = 1.67
$ws.Range("AG6").Value = 800

# Row 11 (Sacramento Republic vs Monterey Bay)
$ws.Range("I11").Value  = 6.6
$ws.Range("L11").Value  = 6.3
$ws.Range("W11").Value  = 6.1
$ws.Range("AH11").Value = 15.5
$ws.Range("AI11").Value = 40
$ws.Range("AM11").Value = 75

# Row 12 (TNS vs Caernarfon)
$ws.Range("G12").Value  = 1.2
$ws.Range("I12").Value  = 12
$ws.Range("J12").Value  = 1.55
$ws.Range("L12").Value  = 8.5
$ws.Range("P12").Value  = 5.7
$ws.Range("R12").Value  = 2.95
$ws.Range("T12").Value  = 3.95
$ws.Range("AH12").Value = 45
$ws.Range("AI12").Value = 120
$ws.Range("AJ12").Value = 37
$ws.Range("AL12").Value = 150
$ws.Range("AN12").Value = 3.4
$ws.Range("AP12").Value = 12.5
$ws.Range("AT12").Value = 3.95
$ws.Range("AU12").Value = 8.25
$ws.Range("AV12").Value = 55
$ws.Range("AW12").Value = 12
$ws.Range("AY12").Value = 40
